$wb = $excel.ActiveWorkbook

# --- Summary sheet: latest_reported_idrc (C), additional_cost_2022 (D), additional_cost_2023 (E), additional_cost_2024 (F) ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 58559136.85168231
$wsSummary.Range("D2").Value = 120825497.2916667
$wsSummary.Range("E2").Value = 90325169.33333333
$wsSummary.Range("F2").Value = 2522152.875

$wsSummary.Range("C3").Value = 238004283.1466135
$wsSummary.Range("D3").Value = 375887809.2166666
$wsSummary.Range("E3").Value = 316196646.2083333
$wsSummary.Range("F3").Value = 12357893.875

$wsSummary.Range("C4").Value = 483317981.6974695
$wsSummary.Range("D4").Value = 731356400.6999999
$wsSummary.Range("E4").Value = 1444882620.6
$wsSummary.Range("F4").Value = 35440798.49999999

$wsSummary.Range("C5").Value = 363211408.4846383
$wsSummary.Range("D5").Value = 1011381128.066666
$wsSummary.Range("E5").Value = 803489231.9333333

$wsSummary.Range("C6").Value = 6480470.109048648
$wsSummary.Range("D6").Value = 1445649651.6
$wsSummary.Range("E6").Value = 1053099530.4
$wsSummary.Range("F6").Value = 32290894.8

$wsSummary.Range("C7").Value = 2567870947.881485
$wsSummary.Range("D7").Value = 8195999177.849999
$wsSummary.Range("E7").Value = 6797462348.85
$wsSummary.Range("F7").Value = 83140754.09999998

$wsSummary.Range("C8").Value = 60858884.68409418
$wsSummary.Range("D8").Value = 456036846.2999999
$wsSummary.Range("E8").Value = 356135839.1999999
$wsSummary.Range("F8").Value = 11517002.1

$wsSummary.Range("C9").Value = 212595686.1506533
$wsSummary.Range("D9").Value = 281514630.4583333
$wsSummary.Range("E9").Value = 220521248.8166667
$wsSummary.Range("F9").Value = 8968522.825000001

$wsSummary.Range("C10").Value = 64843054.67945981
$wsSummary.Range("D10").Value = 401811525.475
$wsSummary.Range("E10").Value = 379755253.4250001

$wsSummary.Range("C11").Value = 1052522864.031498
$wsSummary.Range("D11").Value = 427999756.3666667
$wsSummary.Range("E11").Value = 299719950.0333334

$wsSummary.Range("C12").Value = 1359976138.191323
$wsSummary.Range("D12").Value = 1363847472.5
$wsSummary.Range("E12").Value = 1717765007.5
$wsSummary.Range("F12").Value = 181175310

$wsSummary.Range("C13").Value = 35077534.88405661
$wsSummary.Range("D13").Value = 11800699.21666666
$wsSummary.Range("E13").Value = 11252798.28333333
$wsSummary.Range("F13").Value = 299836.8999999998

$wsSummary.Range("C14").Value = 1276285.374746241
$wsSummary.Range("D14").Value = 130707574.8333333
$wsSummary.Range("E14").Value = 96062394.16666666
$wsSummary.Range("F14").Value = 1335410.2

$wsSummary.Range("C15").Value = 47432854.05086241
$wsSummary.Range("D15").Value = 416711193.8416666
$wsSummary.Range("E15").Value = 528806818.1583333
$wsSummary.Range("F15").Value = 14958868.59999999

$wsSummary.Range("C16").Value = 4314709.876317098
$wsSummary.Range("D16").Value = 6077417.899999999
$wsSummary.Range("E16").Value = 9243425.425000001
$wsSummary.Range("F16").Value = 210550.8749999999

$wsSummary.Range("C17").Value = 510450554.9054651
$wsSummary.Range("D17").Value = 1339734000.916667
$wsSummary.Range("E17").Value = 1098680485.65
$wsSummary.Range("F17").Value = 12097974.53333333

$wsSummary.Range("C18").Value = 224926.0855395237
$wsSummary.Range("D18").Value = 19795
$wsSummary.Range("E18").Value = 18725

$wsSummary.Range("C19").Value = 1966967.651507459
$wsSummary.Range("D19").Value = 162465020.075
$wsSummary.Range("E19").Value = 120113678.6
$wsSummary.Range("F19").Value = 3520167.325

$wsSummary.Range("D20").Value = 533409.7999999999
$wsSummary.Range("E20").Value = 313792.6

$wsSummary.Range("C21").Value = 383424781.3742747
$wsSummary.Range("D21").Value = 1124089390.666667
$wsSummary.Range("E21").Value = 880990474
$wsSummary.Range("F21").Value = 8452362.33333333

$wsSummary.Range("C22").Value = 56366852.26050603
$wsSummary.Range("D22").Value = 424794428.5
$wsSummary.Range("E22").Value = 600741655.1500001
$wsSummary.Range("F22").Value = 49382129.55

$wsSummary.Range("C23").Value = 11052205.68546864
$wsSummary.Range("D23").Value = 71900666.66666667
$wsSummary.Range("E23").Value = 14380133.33333333

$wsSummary.Range("C24").Value = 16303317.58315022
$wsSummary.Range("D24").Value = 1577936519.316667
$wsSummary.Range("E24").Value = 1138536397.558333
$wsSummary.Range("F24").Value = 8680649.825000001

$wsSummary.Range("C25").Value = 10345677.54923651
$wsSummary.Range("D25").Value = 211822587.9
$wsSummary.Range("E25").Value = 150109531.2
$wsSummary.Range("F25").Value = 694221.2999999997

$wsSummary.Range("C26").Value = 1154865.353723471
$wsSummary.Range("D26").Value = 220896471.7499999
$wsSummary.Range("E26").Value = 177892554.45
$wsSummary.Range("F26").Value = 6537642.6

$wsSummary.Range("C27").Value = 2209267.842568598
$wsSummary.Range("D27").Value = 3223243.766666667
$wsSummary.Range("E27").Value = 2515718.958333333
$wsSummary.Range("F27").Value = 22020.675

$wsSummary.Range("C28").Value = 79263094.67844568
$wsSummary.Range("D28").Value = 243333512.1
$wsSummary.Range("E28").Value = 181957010.2
$wsSummary.Range("F28").Value = 6861882.1

$wsSummary.Range("C29").Value = 5072989206.702785
$wsSummary.Range("D29").Value = 520735133.3333333
$wsSummary.Range("E29").Value = 370455866.6666666

# --- Cost per refugee sheet: cost_per_refugee (B) ---
$wsCostPerRefugee = $wb.Worksheets.Item("Cost per refugee")

$wsCostPerRefugee.Range("B2").Value = 2231.5
$wsCostPerRefugee.Range("B3").Value = 9976.1
$wsCostPerRefugee.Range("B4").Value = 12428.1
$wsCostPerRefugee.Range("B5").Value = 23432.8
$wsCostPerRefugee.Range("B6").Value = 5018.4
$wsCostPerRefugee.Range("B7").Value = 14268.6
$wsCostPerRefugee.Range("B8").Value = 19695.6
$wsCostPerRefugee.Range("B9").Value = 2950.9
$wsCostPerRefugee.Range("B10").Value = 16589.9
$wsCostPerRefugee.Range("B11").Value = 6115.6
$wsCostPerRefugee.Range("B12").Value = 16322.1
$wsCostPerRefugee.Range("B13").Value = 1028.6
$wsCostPerRefugee.Range("B14").Value = 6660.4
$wsCostPerRefugee.Range("B15").Value = 12241.3
$wsCostPerRefugee.Range("B16").Value = 5808.3
$wsCostPerRefugee.Range("B17").Value = 13726.9
$wsCostPerRefugee.Range("B18").Value = 21.4
$wsCostPerRefugee.Range("B19").Value = 118
$wsCostPerRefugee.Range("B20").Value = 3737.9
$wsCostPerRefugee.Range("B21").Value = 125.4
$wsCostPerRefugee.Range("B22").Value = 22439.9
$wsCostPerRefugee.Range("B23").Value = 23761.4
$wsCostPerRefugee.Range("B24").Value = 21570.2
$wsCostPerRefugee.Range("B25").Value = 1720.9
$wsCostPerRefugee.Range("B26").Value = 6226.2
$wsCostPerRefugee.Range("B27").Value = 3586.2
$wsCostPerRefugee.Range("B28").Value = 624.7
$wsCostPerRefugee.Range("B29").Value = 8009.2
$wsCostPerRefugee.Range("B30").Value = 10484.6
